$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.473.99"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "2.542.70"
$ws.Range("E3").Value = "  +8.25%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.96"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.77"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  +6.01%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  +12.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.93"
$ws.Range("E10").Value = "  +11.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0828"
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.97"
$ws.Range("E12").Value = "  +11.16%  "
$ws.Range("D13").Value = "2.932.14"
$ws.Range("E13").Value = "  +7.91%  "
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("D15").Value = "2.547.74"
$ws.Range("E15").Value = "  +8.32%  "
$ws.Range("E16").Value = "  +10.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "15.10"
$ws.Range("E17").Value = "  +9.94%  "
$ws.Range("D18").Value = "46.561.21"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.04"
$ws.Range("E19").Value = "  +11.17%  "
$ws.Range("D20").Value = "0.0₃0994"
$ws.Range("E20").Value = "  +3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.61"
$ws.Range("E21").Value = "  +10.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.85"
$ws.Range("E22").Value = "  +5.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.50"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("E24").Value = "  +6.07%  "
$ws.Range("E25").Value = "  +11.73%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.32"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.15"
$ws.Range("E28").Value = "  +15.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.49"
$ws.Range("E29").Value = "  +8.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.78"
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.00"
$ws.Range("E32").Value = "  +10.58%  "
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0844"
$ws.Range("E34").Value = "  +9.41%  "
$ws.Range("E35").Value = "  +21.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.75"
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  +6.09%  "
$ws.Range("E38").Value = "  +4.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.62"
$ws.Range("E39").Value = "  +8.64%  "
$ws.Range("E40").Value = "  +9.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0327"
$ws.Range("E41").Value = "  +9.70%  "
$ws.Range("E42").Value = "  +10.97%  "
$ws.Range("D43").Value = "2.001.85"
$ws.Range("E43").Value = "  +8.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.01"
$ws.Range("E45").Value = "  +3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.61"
$ws.Range("E46").Value = "  +35.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.88"
$ws.Range("E47").Value = "  +2.97%  "
$ws.Range("E48").Value = "  +9.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.97"
$ws.Range("E49").Value = "  +11.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.91"
$ws.Range("E50").Value = "  +11.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.43"
$ws.Range("E51").Value = "  +6.54%  "
